$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.378287315368652
$ws.Range("B1").Value = 1.50287926197052
$ws.Range("C1").Value = 3.860518217086792
$ws.Range("D1").Value = 5.674215316772461
$ws.Range("E1").Value = 1.590700149536133
